# Refresh the cryptos worksheet with the latest scraped price/volume snapshot.
# (GitHub Actions-generated update: prices, 1h volume deltas, and a handful of
# rank re-orderings among the lower-cap coins.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.473.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.07%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.687.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.03%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '682.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.37%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.685.55'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.10%  '

# Row 9
$ws.Range("E9").Value = '  -4.16%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.74%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.52%  '

# Row 13
$ws.Range("E13").Value = '  -4.64%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.03%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.311.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.683.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.502.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.94%  '

# Row 18
$ws.Range("E18").Value = '  -0.72%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.53%  '

# Row 20
$ws.Range("E20").Value = '  -6.47%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '484.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.79%  '

# Row 22
$ws.Range("E22").Value = '  -6.96%  '

# Row 23
$ws.Range("E23").Value = '  -7.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.833.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000131'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.51%  '

# Row 27
$ws.Range("E27").Value = '  -0.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.50'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.87%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.74'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.87%  '

# Row 32
$ws.Range("E32").Value = '  -7.42%  '

# Row 34
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.168'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.62%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.12'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.54%  '

# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.658.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.88%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.50%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0939'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.73%  '

# Row 41
$ws.Range("E41").Value = '  -5.55%  '

# Row 42
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("E43").Value = '  -0.01%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.957'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.18%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '160.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.65%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.40%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.10%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -11.42%  '

# Row 49
$ws.Range("E49").Value = '  -7.30%  '

# Row 50
$ws.Range("E50").Value = '  +0.88%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '394.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.36%  '
